# Update GridConnection's expected validation data (DW) — the equation for
# tower_to_point_of_interconnection_usd_per_kw was fixed, so the validation
# results captured under the GridConnectionCost filter need refreshing, and
# the AutoFilter needs to be switched from ManagementCost back to
# GridConnectionCost so those corrected rows are the ones on display.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39 (ge15_dist_01 / GridConnectionCost / Grid Connection) ----------
# Previously a shared formula (=H39) derived from the old
# tower_to_point_of_interconnection_usd_per_kw equation; now a literal value
# from the corrected equation, formatted like the other recomputed rows
# (the accounting/Currency number format used in rows 136-142).
$ws.Range("G39").NumberFormat = $ws.Range("G136").NumberFormat
$ws.Range("H39").NumberFormat = $ws.Range("H136").NumberFormat
$ws.Range("I39").NumberFormat = $ws.Range("I136").NumberFormat
$ws.Range("G39").Value = 356384.85580801498
$ws.Range("H39").Value = 356384.85580801498
$ws.Range("I39").Value = 237.58990387201001

# --- Row 66 (ge15_dist_05 / GridConnectionCost / Grid Connection) ---------
$ws.Range("G66").NumberFormat = $ws.Range("G136").NumberFormat
$ws.Range("H66").NumberFormat = $ws.Range("H136").NumberFormat
$ws.Range("I66").NumberFormat = $ws.Range("I136").NumberFormat
$ws.Range("G66").Value = 230037.90286025999
$ws.Range("H66").Value = 1150189.5143013
$ws.Range("I66").Value = 153.35860190683999

# --- Row 93 (ge15_dist_10 / GridConnectionCost / Grid Connection) ---------
$ws.Range("G93").NumberFormat = $ws.Range("G136").NumberFormat
$ws.Range("H93").NumberFormat = $ws.Range("H136").NumberFormat
$ws.Range("I93").NumberFormat = $ws.Range("I136").NumberFormat
$ws.Range("G93").Value = 190510.64214379399
$ws.Range("H93").Value = 1905106.4214379401
$ws.Range("I93").Value = 127.00709476252899

# --- Switch the AutoFilter over to GridConnectionCost ----------------------
# This flips which rows are hidden/visible throughout the sheet (the
# ManagementCost rows become hidden, the GridConnectionCost rows become
# visible) to match the new filter selection.
$rng = $ws.Range("A1:I241")
$rng.AutoFilter(4, @("GridConnectionCost"), 7)

# --- Move the active selection to reflect where the user ended up ----------
$ws.Range("H250").Select()

Write-Host "Applied GridConnectionCost validation data fix"
